# Commit swaps the two theme parts in the package:
#   ppt/theme/theme1.xml (the theme used by the slide master / all slides)
#   ppt/theme/theme2.xml (the theme used by the notes master)
# so that the slide master ends up using the "Office" color scheme while
# the notes master ends up using the original "Integral" / "Red Violet"
# color scheme. Both theme parts already share an identical font scheme
# and format scheme (fills/lines/effects), so the only real content
# difference between the two files is their 12-slot theme color scheme
# (and the cosmetic <a:theme>/<a:clrScheme> "name" attributes, which are
# not independently exposed by the object model and simply track the
# color values).
#
# The PowerPoint object model only exposes one live ThemeColorScheme for
# the whole deck (reachable from the slide master, any slide, or the
# notes/handout master alike), and it is backed by theme1.xml -- so that
# is the theme we drive to the new ("Office") palette here.

function ToComRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Color order exposed by ThemeColorScheme.Colors(index):
# 1=dk1 2=lt1 3=dk2 4=lt2 5=accent1 6=accent2 7=accent3 8=accent4
# 9=accent5 10=accent6 11=hlink 12=folHlink
$officeColors = @(
    "000000", "FFFFFF", "44546A", "E7E6E6",
    "5B9BD5", "ED7D31", "A5A5A5", "FFC000",
    "4472C4", "70AD47", "0563C1", "954F72"
)

$p = $ppt.ActivePresentation
$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme

for ($i = 1; $i -le 12; $i++) {
    $colorScheme.Colors($i).RGB = ToComRgb $officeColors[$i - 1]
}
